$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.597.18"
$ws.Range("E2").Value = "  +1.23%  "

$ws.Range("D3").Value = "1.796.61"

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.559"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.26%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.296"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0696"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "

$ws.Range("E11").Value = "  +0.42%  "

$ws.Range("D12").Value = "2.056.49"
$ws.Range("E12").Value = "  +0.86%  "

$ws.Range("D13").Value = "1.795.14"
$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.636"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.29%  "

$ws.Range("D16").Value = "34.576.76"
$ws.Range("E16").Value = "  +1.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0810"
$ws.Range("E19").Value = "  +1.84%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "247.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.32%  "

$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("E23").Value = "  +1.91%  "

$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "167.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.74%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.116"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.29%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.94%  "

$ws.Range("E31").Value = "  +1.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0524"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "

$ws.Range("E33").Value = "  +2.06%  "

$ws.Range("E34").Value = "  +2.75%  "

$ws.Range("D35").Value = "1.426.23"
$ws.Range("E35").Value = "  -1.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.671"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.64%  "

$ws.Range("E38").Value = "  +1.85%  "

$ws.Range("E39").Value = "  +0.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.90%  "

$ws.Range("E41").Value = "  +1.57%  "

$ws.Range("E42").Value = "  +3.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.934"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0528"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.74%  "

$ws.Range("E47").Value = "  +0.62%  "

$ws.Range("D48").Value = "1.955.95"
$ws.Range("E48").Value = "  +0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.24%  "

$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("E51").Value = "  -5.65%  "
